$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.3668248686718916, 1.5291322766718074, -0.9931304002363741, 2.3654466022045453)"
$ws.Range("C2").Value = "JSU(-1.0568984919693878, 1.3418273652357704, 3.467906409217938, 6.473877891100161)"
$ws.Range("D2").Value = "NIG(0.9878716631290815, 0.7865431130628169, 0.686989032978522, 3.120549948133668)"
$ws.Range("E2").Value = "NIG(2.448078256412928, 2.0129116242838063, 1.5787017903524658, 7.0299846736389995)"
